$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Fix row 148 column D (remove stray leading space in shared string) ---
$ws.Range("D148").Value = 'Nils, Manu, Aleks, Til'

# --- Prepare formatting for the new rows 152-161 by copying formats from row 151 ---
$ws.Range("A151:G151").Copy()
$ws.Range("A152:G161").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New CFP rows 152-161 ---
# Row 152
$ws.Range("A152").Value = 152
$ws.Range("B152").Value = '02/2023'
$ws.Range("C152").Value = 'DeFiChainWiki '
$ws.Range("D152").Value = 'Joshua Kummer'
$ws.Range("E152").Value = 'https://www.reddit.com/r/defiblockchain/comments/10l6451/cfp_defichainwiki_9000_dfi/'
$ws.Hyperlinks.Add($ws.Range("E152"), 'https://www.reddit.com/r/defiblockchain/comments/10l6451/cfp_defichainwiki_9000_dfi/')
$ws.Range("F152").Value = 9000
$ws.Range("G152").Value = 'passed'

# Row 153
$ws.Range("A153").Value = 153
$ws.Range("B153").Value = '02/2023'
$ws.Range("C153").Value = ',InTheMarket’’ Podcast'
$ws.Range("D153").Value = 'Patrick (Peddy)'
$ws.Range("E153").Value = 'https://www.reddit.com/r/defiblockchain/comments/1051xfw/cfp_inthemarket_podcast'
$ws.Hyperlinks.Add($ws.Range("E153"), 'https://www.reddit.com/r/defiblockchain/comments/1051xfw/cfp_inthemarket_podcast')
$ws.Range("F153").Value = 9000
$ws.Range("G153").Value = 'passed'

# Row 154
$ws.Range("A154").Value = 154
$ws.Range("B154").Value = '02/2023'
$ws.Range("C154").Value = 'Ongoing Investigation/Case Atomic Swap-dBTC Exploit on DeFiChain - Budget 2023'
$ws.Range("D154").Value = 'Lord Mark'
$ws.Range("E154").Value = 'https://www.reddit.com/r/defiblockchain/comments/10wrrt2/cfp_ongoing_investigationcase_atomic_swapdbtc/'
$ws.Hyperlinks.Add($ws.Range("E154"), 'https://www.reddit.com/r/defiblockchain/comments/10wrrt2/cfp_ongoing_investigationcase_atomic_swapdbtc/')
$ws.Range("F154").Value = 200000
$ws.Range("G154").Value = 'passed'

# Row 155
$ws.Range("A155").Value = 155
$ws.Range("B155").Value = '02/2023'
$ws.Range("C155").Value = 'Sponsorship of the National Fighting Championship (NFC)'
$ws.Range("D155").Value = 'DeFiChain Epic, DeFiChain Accelerator'
$ws.Range("E155").Value = 'https://github.com/DeFiCh/dfips/issues/246'
$ws.Hyperlinks.Add($ws.Range("E155"), 'https://github.com/DeFiCh/dfips/issues/246')
$ws.Range("F155").Value = 500000
$ws.Range("G155").Value = 'passed'

# Row 156
$ws.Range("A156").Value = 156
$ws.Range("B156").Value = '02/2023'
$ws.Range("C156").Value = 'defichain-trader.com'
$ws.Range("D156").Value = 'Ruben'
$ws.Range("E156").Value = 'https://www.reddit.com/r/defiblockchain/comments/10l35aj/cfp_defichaintradercom/'
$ws.Hyperlinks.Add($ws.Range("E156"), 'https://www.reddit.com/r/defiblockchain/comments/10l35aj/cfp_defichaintradercom/')
$ws.Range("F156").Value = 110000
$ws.Range("G156").Value = 'declined'

# Row 157
$ws.Range("A157").Value = 157
$ws.Range("B157").Value = '02/2023'
$ws.Range("C157").Value = 'Portfolio Optimisation with Modern Portfolio Theory'
$ws.Range("D157").Value = 'Lukas'
$ws.Range("E157").Value = 'https://www.reddit.com/r/defiblockchain/comments/10x7jhn/cfp_portfolio_optimisation_with_modern_portfolio/'
$ws.Hyperlinks.Add($ws.Range("E157"), 'https://www.reddit.com/r/defiblockchain/comments/10x7jhn/cfp_portfolio_optimisation_with_modern_portfolio/')
$ws.Range("F157").Value = 20000
$ws.Range("G157").Value = 'declined'

# Row 158
$ws.Range("A158").Value = 158
$ws.Range("B158").Value = '04/2023'
$ws.Range("C158").Value = 'defichain-trader.com (Maintenance)'
$ws.Range("D158").Value = 'Ruben'
$ws.Range("E158").Value = 'https://github.com/DeFiCh/dfips/issues/258'
$ws.Hyperlinks.Add($ws.Range("E158"), 'https://github.com/DeFiCh/dfips/issues/258')
$ws.Range("F158").Value = 5400
$ws.Range("G158").Value = 'passed'

# Row 159
$ws.Range("A159").Value = 159
$ws.Range("B159").Value = '04/2023'
$ws.Range("C159").Value = 'Appreciation of the work done by Kügi in the community in the last months'
$ws.Range("D159").Value = 'Phigo'
$ws.Range("E159").Value = 'https://www.reddit.com/r/defiblockchain/comments/11fj7i5/cfp_appreciation_of_the_work_done_by_k%C3%BCgi_in_the/'
$ws.Hyperlinks.Add($ws.Range("E159"), 'https://www.reddit.com/r/defiblockchain/comments/11fj7i5/cfp_appreciation_of_the_work_done_by_k%C3%BCgi_in_the/')
$ws.Range("F159").Value = 20000
$ws.Range("G159").Value = 'passed'

# Row 160
$ws.Range("A160").Value = 160
$ws.Range("B160").Value = '04/2023'
$ws.Range("C160").Value = 'mydefichain-Ocean-2023'
$ws.Range("D160").Value = 'Bernd Mack and Andreas Lentz with mydefichain'
$ws.Range("E160").Value = 'https://www.reddit.com/r/defiblockchain/comments/11g01xa/cfp_mydefichain_ocean_funding_2023_6850_dfi_per/'
$ws.Hyperlinks.Add($ws.Range("E160"), 'https://www.reddit.com/r/defiblockchain/comments/11g01xa/cfp_mydefichain_ocean_funding_2023_6850_dfi_per/')
$ws.Range("F160").Value = 6850
$ws.Range("G160").Value = 'passed'

# Row 161
$ws.Range("A161").Value = 161
$ws.Range("B161").Value = '04/2023'
$ws.Range("C161").Value = 'Continue DFI rewards on BSC DFI-BNB pool on BSC/ACSI.finance'
$ws.Range("D161").Value = 'Michael for/with ACryptoS Team'
$ws.Range("E161").Value = 'https://github.com/DeFiCh/dfips/issues/256'
$ws.Hyperlinks.Add($ws.Range("E161"), 'https://github.com/DeFiCh/dfips/issues/256')
$ws.Range("F161").Value = 20000
$ws.Range("G161").Value = 'declined'

# --- Column widths (best effort; engine quantizes to 1/6-character steps) ---
$ws.Columns.Item(1).ColumnWidth = 5.166666666666667
$ws.Columns.Item(3).ColumnWidth = 91.87760416666667

# --- View state: zoom, scroll position, selection ---
$win = $excel.ActiveWindow
$win.Zoom = 115
$win.ScrollRow = 123
$win.ScrollColumn = 1
$ws.Range("C168").Select()
